$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 321.27274
$ws.Range("I28").Value = 326
$ws.Range("K28").Value = 326
$ws.Range("M28").Value = 159
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H62").Value = 4867.923
$ws.Range("I62").Value = 2910.375
$ws.Range("K62").Value = 2910.375
$ws.Range("M62").Value = -2286.375
$ws.Range("H65").Value = 4867.923
$ws.Range("I65").Value = 2910.375
$ws.Range("K65").Value = 14551.875
$ws.Range("M65").Value = -11431.875
$ws.Range("H86").Value = 3866.125
$ws.Range("I86").Value = 2284.8572
$ws.Range("K86").Value = 2284.8572
$ws.Range("M86").Value = -1161.8572
$ws.Range("H89").Value = 3866.125
$ws.Range("I89").Value = 2284.8572
$ws.Range("K89").Value = 11424.286
$ws.Range("M89").Value = -5808.286
$ws.Range("H98").Value = 1290
$ws.Range("I98").Value = 1076.25
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 1076.25
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 421.75
$ws.Range("N98").Value = -5996
$ws.Range("H106").Value = 22371.75
$ws.Range("I106").Value = 27273.445
$ws.Range("K106").Value = 27273.445
$ws.Range("M106").Value = -26642.445
$ws.Range("H107").Value = 464.72726
$ws.Range("I107").Value = 505.125
$ws.Range("J107").Value = 357
$ws.Range("K107").Value = 505.125
$ws.Range("L107").Value = 357
$ws.Range("M107").Value = 1414.875
$ws.Range("N107").Value = -4197
$ws.Range("H116").Value = 7470
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 7470
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 7470
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -14354
$ws.Range("H122").Value = 1290
$ws.Range("I122").Value = 1076.25
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3228.75
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -778.75
$ws.Range("N122").Value = -13900
$ws.Range("H125").Value = 1497.3334
$ws.Range("I125").Value = 1497.3334
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 13476.0006
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -11016.0006
$ws.Range("N125").ClearContents()
$ws.Range("H137").Value = 3656.5715
$ws.Range("I137").Value = 1000
$ws.Range("K137").Value = 3000
$ws.Range("M137").Value = -450
$ws.Range("H138").Value = 4207.597
$ws.Range("J138").Value = 4515.932
$ws.Range("L138").Value = 13547.796
$ws.Range("N138").Value = -23827.796
$ws.Range("H141").Value = 5138.1816
$ws.Range("I141").Value = 5002.222
$ws.Range("J141").Value = 5750
$ws.Range("K141").Value = 15006.666
$ws.Range("L141").Value = 17250
$ws.Range("M141").Value = -9826.665999999999
$ws.Range("N141").Value = -27610

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1518.6
$ws.Range("I61").Value = 1523.5
$ws.Range("J61").Value = 1499
$ws.Range("K61").Value = 1523.5
$ws.Range("L61").Value = 1499
$ws.Range("M61").Value = -1311.5
$ws.Range("N61").Value = -1923
$ws.Range("H74").Value = 4840.8184
$ws.Range("I74").Value = 1583.1666
$ws.Range("K74").Value = 1583.1666
$ws.Range("M74").Value = -709.1666
$ws.Range("H77").Value = 4840.8184
$ws.Range("I77").Value = 1583.1666
$ws.Range("K77").Value = 7915.833000000001
$ws.Range("M77").Value = -3547.833000000001
$ws.Range("H110").Value = 3138.1538
$ws.Range("I110").Value = 3536
$ws.Range("K110").Value = 3536
$ws.Range("M110").Value = -1491
$ws.Range("H122").Value = 836761.75
$ws.Range("I122").Value = 1114015.6
$ws.Range("K122").Value = 3342046.8
$ws.Range("M122").Value = -3339596.8
$ws.Range("H132").Value = 1341.3096
$ws.Range("I132").Value = 1404.0264
$ws.Range("J132").Value = 745.5
$ws.Range("K132").Value = 4212.0792
$ws.Range("L132").Value = 2236.5
$ws.Range("M132").Value = -1682.0792
$ws.Range("N132").Value = -7296.5
$ws.Range("H136").Value = 1518.6
$ws.Range("I136").Value = 1523.5
$ws.Range("J136").Value = 1499
$ws.Range("K136").Value = 4570.5
$ws.Range("L136").Value = 4497
$ws.Range("M136").Value = -2020.5
$ws.Range("N136").Value = -9597

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 999.3333
$ws.Range("I94").Value = 999.3333
$ws.Range("K94").Value = 999.3333
$ws.Range("M94").Value = -548.3333
$ws.Range("H105").Value = 3046.9285
$ws.Range("I105").Value = 2324.9355
$ws.Range("K105").Value = 2324.9355
$ws.Range("M105").Value = -577.9355
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H134").Value = 3373.1155
$ws.Range("I134").Value = 3425.5417
$ws.Range("J134").Value = 2744
$ws.Range("K134").Value = 10276.6251
$ws.Range("L134").Value = 8232
$ws.Range("M134").Value = -7741.625100000001
$ws.Range("N134").Value = -13302

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 60000
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -62246
$ws.Range("H83").Value = 60000
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -191232
$ws.Range("H94").Value = 1998
$ws.Range("I94").Value = 1998
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1998
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1547
$ws.Range("N94").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 638.0714
$ws.Range("J5").Value = 862.2857
$ws.Range("L5").Value = 2586.8571
$ws.Range("N5").Value = -2810.8571
$ws.Range("H33").Value = 1359.625
$ws.Range("I33").Value = 3406.6667
$ws.Range("J33").Value = 131.4
$ws.Range("K33").Value = 20440.0002
$ws.Range("L33").Value = 788.4000000000001
$ws.Range("M33").Value = -20157.0002
$ws.Range("N33").Value = -1354.4
$ws.Range("H62").Value = 4574.875
$ws.Range("J62").Value = 4200
$ws.Range("L62").Value = 12600
$ws.Range("N62").Value = -13972
$ws.Range("H64").Value = 1215.6
$ws.Range("I64").Value = 583
$ws.Range("J64").Value = 1637.3334
$ws.Range("K64").Value = 1749
$ws.Range("L64").Value = 4912.0002
$ws.Range("M64").Value = -1479
$ws.Range("N64").Value = -5452.0002
$ws.Range("H65").Value = 4574.875
$ws.Range("J65").Value = 4200
$ws.Range("L65").Value = 37800
$ws.Range("N65").Value = -44664
$ws.Range("H67").Value = 1215.6
$ws.Range("I67").Value = 583
$ws.Range("J67").Value = 1637.3334
$ws.Range("K67").Value = 1749
$ws.Range("L67").Value = 4912.0002
$ws.Range("M67").Value = -813
$ws.Range("N67").Value = -6784.0002
$ws.Range("H80").Value = 3982.7778
$ws.Range("I80").Value = 2722.75
$ws.Range("J80").Value = 4990.8
$ws.Range("K80").Value = 8168.25
$ws.Range("L80").Value = 14972.4
$ws.Range("M80").Value = -7232.25
$ws.Range("N80").Value = -16844.4
$ws.Range("H83").Value = 3982.7778
$ws.Range("I83").Value = 2722.75
$ws.Range("J83").Value = 4990.8
$ws.Range("K83").Value = 24504.75
$ws.Range("L83").Value = 44917.2
$ws.Range("M83").Value = -19824.75
$ws.Range("N83").Value = -54277.2
$ws.Range("H86").Value = 341
$ws.Range("I86").Value = 340
$ws.Range("J86").Value = 342
$ws.Range("K86").Value = 1020
$ws.Range("L86").Value = 1026
$ws.Range("M86").Value = 166
$ws.Range("N86").Value = -3398
$ws.Range("H89").Value = 341
$ws.Range("I89").Value = 340
$ws.Range("J89").Value = 342
$ws.Range("K89").Value = 3060
$ws.Range("L89").Value = 3078
$ws.Range("M89").Value = 2868
$ws.Range("N89").Value = -14934
$ws.Range("H107").Value = 1136.75
$ws.Range("I107").Value = 589.5
$ws.Range("K107").Value = 1768.5
$ws.Range("M107").Value = 151.5
$ws.Range("H113").Value = 886.46875
$ws.Range("I113").Value = 1529.6
$ws.Range("J113").Value = 767.37036
$ws.Range("K113").Value = 4588.799999999999
$ws.Range("L113").Value = 2302.11108
$ws.Range("M113").Value = -2418.799999999999
$ws.Range("N113").Value = -6642.111080000001
$ws.Range("H120").Value = 11583.793
$ws.Range("J120").Value = 12107.407
$ws.Range("L120").Value = 36322.221
$ws.Range("N120").Value = -45998.221
$ws.Range("H126").Value = 400
$ws.Range("I126").Value = 400
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1200
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 3740
$ws.Range("N126").ClearContents()
$ws.Range("H135").Value = 638.0714
$ws.Range("J135").Value = 862.2857
$ws.Range("L135").Value = 7760.571300000001
$ws.Range("N135").Value = -12830.5713

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 30000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 30000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 30000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -30224
$ws.Range("H8").Value = 30000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 30000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 30000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -30278
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H29").Value = 17500
$ws.Range("I29").Value = 15400
$ws.Range("K29").Value = 15400
$ws.Range("M29").Value = -15110
$ws.Range("H35").Value = 5500000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H97").Value = 1980.5
$ws.Range("I97").Value = 1925.1111
$ws.Range("K97").Value = 1925.1111
$ws.Range("M97").Value = -1429.1111
$ws.Range("H107").Value = 953.9545000000001
$ws.Range("J107").Value = 965.2632
$ws.Range("L107").Value = 965.2632
$ws.Range("N107").Value = -4805.2632
$ws.Range("H113").Value = 2457.111
$ws.Range("I113").Value = 1185.6666
$ws.Range("K113").Value = 1185.6666
$ws.Range("M113").Value = 984.3334
$ws.Range("H126").Value = 5112.3335
$ws.Range("I126").Value = 3006
$ws.Range("K126").Value = 9018
$ws.Range("M126").Value = -6548
$ws.Range("H132").Value = 3306.6667
$ws.Range("I132").Value = 2976.4736
$ws.Range("K132").Value = 8929.4208
$ws.Range("M132").Value = -6399.4208

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2572.75
$ws.Range("I7").Value = 2572.75
$ws.Range("K7").Value = 2572.75
$ws.Range("M7").Value = -2460.75
$ws.Range("H16").Value = 3660.44
$ws.Range("I16").Value = 3566.95
$ws.Range("K16").Value = 3566.95
$ws.Range("M16").Value = -3396.95
$ws.Range("H20").Value = 16253
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 16253
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 16253
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -16705
$ws.Range("H21").Value = 1849.75
$ws.Range("I21").Value = 306
$ws.Range("J21").Value = 2364.3333
$ws.Range("K21").Value = 306
$ws.Range("L21").Value = 2364.3333
$ws.Range("M21").Value = -132
$ws.Range("N21").Value = -2712.3333
$ws.Range("H61").Value = 2440.4707
$ws.Range("I61").Value = 1883.5385
$ws.Range("K61").Value = 1883.5385
$ws.Range("M61").Value = -1681.5385
$ws.Range("H80").Value = 34999
$ws.Range("J80").Value = 34999
$ws.Range("L80").Value = 34999
$ws.Range("N80").Value = -37245
$ws.Range("H83").Value = 34999
$ws.Range("J83").Value = 34999
$ws.Range("L83").Value = 104997
$ws.Range("N83").Value = -116229
$ws.Range("H113").Value = 2440.4707
$ws.Range("I113").Value = 1883.5385
$ws.Range("K113").Value = 1883.5385
$ws.Range("M113").Value = 286.4614999999999
$ws.Range("H122").Value = 1000
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 2572.75
$ws.Range("I126").Value = 2572.75
$ws.Range("K126").Value = 7718.25
$ws.Range("M126").Value = -5248.25
$ws.Range("H132").Value = 5125.6313
$ws.Range("I132").Value = 4733.64
$ws.Range("K132").Value = 14200.92
$ws.Range("M132").Value = -11670.92

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1237.5883
$ws.Range("I122").Value = 1237.5883
$ws.Range("K122").Value = 3712.7649
$ws.Range("M122").Value = -1262.7649
$ws.Range("H136").Value = 51575.75
$ws.Range("I136").Value = 970.86664
$ws.Range("K136").Value = 2912.59992
$ws.Range("M136").Value = -362.5999199999997
